# Auto-generated edit script applying the scheduled-runner market-price update
# to the Cerberus_Profits leve-profit tables (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 5642.5527
$ws.Range("J74").Value = 5710.88
$ws.Range("L74").Value = 5710.88
$ws.Range("N74").Value = -7582.88
$ws.Range("H77").Value = 5642.5527
$ws.Range("J77").Value = 5710.88
$ws.Range("L77").Value = 28554.4
$ws.Range("N77").Value = -37914.4
$ws.Range("H112").Value = 5294.5713
$ws.Range("J112").Value = 5837.087
$ws.Range("L112").Value = 17511.261
$ws.Range("N112").Value = -19727.261
$ws.Range("H129").Value = 1566
$ws.Range("J129").Value = 4100
$ws.Range("L129").Value = 12300
$ws.Range("N129").Value = -22300
$ws.Range("H135").Value = 1975.125
$ws.Range("I135").Value = 1800.1818
$ws.Range("K135").Value = 16201.6362
$ws.Range("M135").Value = -13666.6362
$ws.Range("H138").Value = 3470.8147
$ws.Range("I138").Value = 4813.5884
$ws.Range("J138").Value = 2853.8647
$ws.Range("K138").Value = 14440.7652
$ws.Range("L138").Value = 8561.5941
$ws.Range("M138").Value = -9300.765199999998
$ws.Range("N138").Value = -18841.5941

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2534.2666
$ws.Range("I45").Value = 1961.7
$ws.Range("J45").Value = 3679.4
$ws.Range("K45").Value = 1961.7
$ws.Range("L45").Value = 3679.4
$ws.Range("M45").Value = -1584.7
$ws.Range("N45").Value = -4433.4

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 9999
$ws.Range("I75").Value = 9999
$ws.Range("K75").Value = 9999
$ws.Range("M75").Value = -9063
$ws.Range("H78").Value = 9999
$ws.Range("I78").Value = 9999
$ws.Range("K78").Value = 29997
$ws.Range("M78").Value = -25317
$ws.Range("H99").Value = 4952.6
$ws.Range("I99").Value = 4899.5
$ws.Range("K99").Value = 4899.5
$ws.Range("M99").Value = -3401.5

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H37").Value = 1928.5714
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("H52").Value = 76900
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H58").Value = 1370.8948
$ws.Range("I58").Value = 980.6923
$ws.Range("J58").Value = 2216.3333
$ws.Range("K58").Value = 980.6923
$ws.Range("L58").Value = 2216.3333
$ws.Range("M58").Value = -777.6923
$ws.Range("N58").Value = -2622.3333
$ws.Range("H136").Value = 1370.8948
$ws.Range("I136").Value = 980.6923
$ws.Range("J136").Value = 2216.3333
$ws.Range("K136").Value = 2942.0769
$ws.Range("L136").Value = 6648.999899999999
$ws.Range("M136").Value = -392.0769
$ws.Range("N136").Value = -11748.9999

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 1731.1666
$ws.Range("J75").Value = 1718.5
$ws.Range("L75").Value = 5155.5
$ws.Range("N75").Value = -7151.5
$ws.Range("H78").Value = 1731.1666
$ws.Range("J78").Value = 1718.5
$ws.Range("L78").Value = 15466.5
$ws.Range("N78").Value = -25450.5
$ws.Range("H117").Value = 66670120
$ws.Range("I117").Value = 111112620
$ws.Range("J117").Value = 47623336
$ws.Range("K117").Value = 333337860
$ws.Range("L117").Value = 142870008
$ws.Range("M117").Value = -333334418
$ws.Range("N117").Value = -142876892
$ws.Range("H118").Value = 3574
$ws.Range("I118").Value = 3574
$ws.Range("K118").Value = 10722
$ws.Range("M118").Value = -9479
$ws.Range("H122").Value = 1082.5555
$ws.Range("I122").Value = 248.6
$ws.Range("J122").Value = 2125
$ws.Range("K122").Value = 2237.4
$ws.Range("L122").Value = 19125
$ws.Range("M122").Value = 212.5999999999999
$ws.Range("N122").Value = -24025
$ws.Range("H134").Value = 13899.523
$ws.Range("I134").Value = 1380
$ws.Range("J134").Value = 17811.875
$ws.Range("K134").Value = 4140
$ws.Range("L134").Value = 53435.625
$ws.Range("M134").Value = 930
$ws.Range("N134").Value = -63575.625
$ws.Range("H136").Value = 5754
$ws.Range("I136").Value = 3114.8
$ws.Range("J136").Value = 18950
$ws.Range("K136").Value = 9344.400000000001
$ws.Range("L136").Value = 56850
$ws.Range("M136").Value = -4244.400000000001
$ws.Range("N136").Value = -67050
$ws.Range("H137").Value = 17149.25
$ws.Range("I137").Value = 18148
$ws.Range("J137").Value = 16550
$ws.Range("K137").Value = 54444
$ws.Range("L137").Value = 49650
$ws.Range("M137").Value = -49344
$ws.Range("N137").Value = -59850

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 316.07693
$ws.Range("I107").Value = 175.75
$ws.Range("K107").Value = 175.75
$ws.Range("M107").Value = 1744.25
$ws.Range("H132").Value = 2189.7896
$ws.Range("I132").Value = 1440.25
$ws.Range("K132").Value = 4320.75
$ws.Range("M132").Value = -1790.75

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5010.72
$ws.Range("J122").Value = 5681.8887
$ws.Range("L122").Value = 17045.6661
$ws.Range("N122").Value = -21945.6661
$ws.Range("H132").Value = 2139.8572
$ws.Range("I132").Value = 1635.2941
$ws.Range("K132").Value = 4905.8823
$ws.Range("M132").Value = -2375.8823

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 52986.223
$ws.Range("J46").Value = 52986.223
$ws.Range("L46").Value = 52986.223
$ws.Range("N46").Value = -53448.223
$ws.Range("H62").Value = 5221.3335
$ws.Range("I62").Value = 3332
$ws.Range("K62").Value = 3332
$ws.Range("M62").Value = -2708
$ws.Range("H65").Value = 5221.3335
$ws.Range("I65").Value = 3332
$ws.Range("K65").Value = 16660
$ws.Range("M65").Value = -13540
$ws.Range("H101").Value = 69999
$ws.Range("J101").Value = 69999
$ws.Range("L101").Value = 69999
$ws.Range("N101").Value = -76489
$ws.Range("H107").Value = 1220.75
$ws.Range("I107").Value = 441.75
$ws.Range("K107").Value = 1325.25
$ws.Range("M107").Value = 594.75
$ws.Range("H117").Value = 44996
$ws.Range("J117").Value = 44996
$ws.Range("L117").Value = 44996
$ws.Range("N117").Value = -54174
$ws.Range("H134").Value = 52986.223
$ws.Range("J134").Value = 52986.223
$ws.Range("L134").Value = 158958.669
$ws.Range("N134").Value = -164028.669
$ws.Range("H136").Value = 4824.385
$ws.Range("I136").Value = 4830.5
$ws.Range("J136").Value = 4814.6
$ws.Range("K136").Value = 14491.5
$ws.Range("L136").Value = 14443.8
$ws.Range("M136").Value = -11941.5
$ws.Range("N136").Value = -19543.8

